$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-16 16:59:43"

$wsZhCn.Range("H2").Value = "2016-08-16 16:59:38"
$wsZhCn.Range("K2").Value = "2016-08-16 16:59:55"

$wsDeDe.Range("H2").Value = "2016-08-16 17:00:10"
